# Add a "Price" column (N) with a header and per-row stock prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column, styled like the other plain (unshaded) cells.
$ws.Range("N1").Value = "Price"

# Price values for rows 7 through 44 (rows 2-6 stay blank in this column).
$prices = @{
    7  = 31.28
    8  = 32.06
    9  = 30.85
    10 = 41.02
    11 = 38.03
    12 = 38.79
    13 = 45.36
    14 = 41.91
    15 = 35
    16 = 44.33
    17 = 31.25
    18 = 27.45
    19 = 18.420000000000002
    20 = 17.14
    21 = 16.64
    22 = 9.1
    23 = 14.51
    24 = 9.42
    25 = 30.87
    26 = 21.36
    27 = 24.42
    28 = 32.909999999999997
    29 = 32.82
    30 = 37.450000000000003
    31 = 44.38
    32 = 41.45
    33 = 48.12
    34 = 41.53
    35 = 50.01
    36 = 47.56
    37 = 47.56
    38 = 41.53
    39 = 50.01
    40 = 47.56
    41 = 37.47
    42 = 46.75
    43 = 52.67
    44 = 68.37
}

foreach ($row in $prices.Keys) {
    $cell = $ws.Range("N" + $row)
    $cell.Value = $prices[$row]
    $cell.NumberFormat = "#,##0.00;[Red](#,##0.00)"
}

# Restore the selection left behind in the source file.
$ws.Range("P21").Select()
